# Update planning, themes, SRS V1.1
#
# The Sprint 3 (rows 14-18) and Sprint 4 (rows 20-24) task blocks on
# "Sheet2" still had placeholder Module / Module-Desc / Status detail
# text copied over from the Sprint 2 block. Clear that placeholder
# detail content (columns D, E and J) for those two blocks, leaving
# just the block headers / task-name cells intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$cellsToClear = @(
    "D14", "E14", "J14",
    "J15",
    "D16", "E16", "J16",
    "J17",
    "D18", "E18", "J18",
    "D20", "E20", "J20",
    "J21",
    "D22", "E22", "J22",
    "J23",
    "D24", "E24", "J24"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value2 = ""
}

# Leave the sheet scrolled/zoomed where the author ended up after
# reviewing the cleared Sprint 3/4 rows.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 101
$ws.Range("K11").Select() | Out-Null
